# Update "想去人数" (interest count) figures in column F across the
# "展览", "演出" and "全部类型" sheets, reflecting newly generated stats.

$wb = $excel.ActiveWorkbook

$updates = @{}

$updates["展览"] = @{
    "F4"  = 112
    "F5"  = 1724
    "F6"  = 3302
    "F7"  = 971
    "F8"  = 2145
    "F9"  = 2062
    "F10" = 1073
    "F11" = 575
    "F13" = 1644
    "F14" = 364
    "F16" = 27
    "F18" = 149
    "F19" = 1518
    "F20" = 574
    "F21" = 674
    "F22" = 561
    "F23" = 12015
    "F24" = 12025
    "F25" = 886
    "F28" = 7
    "F29" = 298
    "F30" = 1885
    "F31" = 176
    "F32" = 512
}

$updates["演出"] = @{
    "F7" = 10
}

$updates["全部类型"] = @{
    "F6"  = 112
    "F7"  = 1724
    "F8"  = 3302
    "F9"  = 971
    "F10" = 2145
    "F11" = 2062
    "F12" = 1073
    "F13" = 575
    "F15" = 1644
    "F16" = 364
    "F18" = 27
    "F22" = 149
    "F23" = 1518
    "F24" = 574
    "F25" = 674
    "F26" = 561
    "F27" = 12015
    "F28" = 12025
    "F29" = 886
    "F32" = 7
    "F33" = 298
    "F34" = 1885
    "F37" = 176
    "F38" = 512
    "F39" = 10
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}

$wb.Save()
